$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing weight values for the InceptionNet / Fixed row (row 7)
$ws.Range("D7").Value = 0.117326826
$ws.Range("E7").Value = 2.5623705000000001
$ws.Range("F7").Value = 0.091901064000000005
$ws.Range("G7").Value = 1.9458221
$ws.Range("H7").Value = 0.079053940000000003
$ws.Range("I7").Value = 1.6741633

# Update the active selection to reflect where the user left off
$ws.Range("G7").Select()
